$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price column cells being updated so Excel does not
# auto-convert numeric-looking strings (e.g. "1.00", "0.150") into numbers,
# matching the source data which stores these as plain text.
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D10", "D13", "D14", "D15", "D16", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply updated cell values as scraped by the latest GitHub Actions run.
$ws.Range("D2").Value = "69.817.28"
$ws.Range("E2").Value = "  +4.21%  "
$ws.Range("D3").Value = "3.635.11"
$ws.Range("E3").Value = "  +3.31%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "633.24"
$ws.Range("E5").Value = "  +3.89%  "
$ws.Range("D6").Value = "159.85"
$ws.Range("E6").Value = "  +5.09%  "
$ws.Range("D7").Value = "3.632.78"
$ws.Range("E7").Value = "  +3.30%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +2.82%  "
$ws.Range("D10").Value = "0.150"
$ws.Range("E10").Value = "  +7.56%  "
$ws.Range("E11").Value = "  +5.77%  "
$ws.Range("E12").Value = "  +3.80%  "
$ws.Range("D13").Value = "0.0000231"
$ws.Range("E13").Value = "  +5.70%  "
$ws.Range("D14").Value = "33.45"
$ws.Range("E14").Value = "  +6.14%  "
$ws.Range("D15").Value = "4.252.60"
$ws.Range("E15").Value = "  +3.37%  "
$ws.Range("D16").Value = "3.637.01"
$ws.Range("E16").Value = "  +3.28%  "
$ws.Range("D17").Value = "69.717.46"
$ws.Range("E17").Value = "  +4.23%  "
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").Value = "6.68"
$ws.Range("E19").Value = "  +6.46%  "
$ws.Range("D20").Value = "16.07"
$ws.Range("E20").Value = "  +4.91%  "
$ws.Range("D21").Value = "10.18"
$ws.Range("E21").Value = "  +11.53%  "
$ws.Range("D22").Value = "465.65"
$ws.Range("E22").Value = "  +5.21%  "
$ws.Range("D23").Value = "0.645"
$ws.Range("E23").Value = "  +2.67%  "
$ws.Range("D24").Value = "78.87"
$ws.Range("E24").Value = "  +1.50%  "
$ws.Range("D25").Value = "0.0000137"
$ws.Range("D26").Value = "10.74"
$ws.Range("E26").Value = "  +5.12%  "
$ws.Range("D27").Value = "3.779.17"
$ws.Range("E27").Value = "  +3.24%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("E29").Value = "  +14.37%  "
$ws.Range("D30").Value = "2.66"
$ws.Range("E30").Value = "  +5.24%  "
$ws.Range("D31").Value = "1.73"
$ws.Range("E31").Value = "  +5.24%  "
$ws.Range("D32").Value = "0.178"
$ws.Range("E32").Value = "  +12.71%  "
$ws.Range("D33").Value = "6.65"
$ws.Range("E33").Value = "  +8.63%  "
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  -0.22%  "
$ws.Range("D35").Value = "1.98"
$ws.Range("E35").Value = "  +6.51%  "
$ws.Range("D36").Value = "26.64"
$ws.Range("E36").Value = "  +3.87%  "
$ws.Range("D37").Value = "3.632.24"
$ws.Range("E37").Value = "  +3.40%  "
$ws.Range("D38").Value = "8.50"
$ws.Range("E38").Value = "  +6.22%  "
$ws.Range("E39").Value = "  +14.78%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").Value = "0.0933"
$ws.Range("E41").Value = "  +8.71%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").Value = "177.50"
$ws.Range("E43").Value = "  +2.21%  "
$ws.Range("D44").Value = "5.65"
$ws.Range("E44").Value = "  +1.79%  "
$ws.Range("D45").Value = "31.83"
$ws.Range("E45").Value = "  +17.93%  "
$ws.Range("D46").Value = "0.917"
$ws.Range("E46").Value = "  +3.23%  "
$ws.Range("D47").Value = "1.39"
$ws.Range("E47").Value = "  +13.91%  "
$ws.Range("D48").Value = "2.81"
$ws.Range("E48").Value = "  +11.14%  "
$ws.Range("D49").Value = "46.52"
$ws.Range("E49").Value = "  +2.96%  "
$ws.Range("E50").Value = "  +3.78%  "
$ws.Range("D51").Value = "0.269"
$ws.Range("E51").Value = "  +8.66%  "
